$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Sat Feb 17 22:44:45 EST 2024"
$ws.Range("B3").Value = "Sat Feb 17 22:44:57 EST 2024"
$ws.Range("B4").Value = "Sat Feb 17 22:45:09 EST 2024"
$ws.Range("B5").Value = "Sat Feb 17 22:45:21 EST 2024"
$ws.Range("B6").Value = "Sat Feb 17 22:45:33 EST 2024"
$ws.Range("B7").Value = "Sat Feb 17 22:45:45 EST 2024"
$ws.Range("B8").Value = "Sat Feb 17 22:45:57 EST 2024"
$ws.Range("B9").Value = "Sat Feb 17 22:46:09 EST 2024"
$ws.Range("B10").Value = "Sat Feb 17 22:46:21 EST 2024"
$ws.Range("B11").Value = "Sat Feb 17 22:46:34 EST 2024"
$ws.Range("B12").Value = "Sat Feb 17 22:46:46 EST 2024"
$ws.Range("B13").Value = "Sat Feb 17 22:46:58 EST 2024"
$ws.Range("B14").Value = "Sat Feb 17 22:47:10 EST 2024"
$ws.Range("B15").Value = "Sat Feb 17 22:47:22 EST 2024"
$ws.Range("B16").Value = "Sat Feb 17 22:47:33 EST 2024"
$ws.Range("B17").Value = "Sat Feb 17 22:47:45 EST 2024"
$ws.Range("B18").Value = "Sat Feb 17 22:47:57 EST 2024"
$ws.Range("B19").Value = "Sat Feb 17 22:48:10 EST 2024"
$ws.Range("B20").Value = "Sat Feb 17 22:48:22 EST 2024"
$ws.Range("B28").Value = "Sat Feb 17 22:48:34 EST 2024"
$ws.Range("B29").Value = "Sat Feb 17 22:48:46 EST 2024"
$ws.Range("B30").Value = "Sat Feb 17 22:48:58 EST 2024"
$ws.Range("B31").Value = "Sat Feb 17 22:49:10 EST 2024"
$ws.Range("B32").Value = "Sat Feb 17 22:49:23 EST 2024"
$ws.Range("B33").Value = "Sat Feb 17 22:49:35 EST 2024"
$ws.Range("B34").Value = "Sat Feb 17 22:49:47 EST 2024"
$ws.Range("B35").Value = "Sat Feb 17 22:49:59 EST 2024"
$ws.Range("B36").Value = "Sat Feb 17 22:50:11 EST 2024"
$ws.Range("B37").Value = "Sat Feb 17 22:50:23 EST 2024"
$ws.Range("B38").Value = "Sat Feb 17 22:50:35 EST 2024"
$ws.Range("B39").Value = "Sat Feb 17 22:50:47 EST 2024"
$ws.Range("B40").Value = "Sat Feb 17 22:50:59 EST 2024"
$ws.Range("B41").Value = "Sat Feb 17 22:51:11 EST 2024"
$ws.Range("B42").Value = "Sat Feb 17 22:51:23 EST 2024"
$ws.Range("B43").Value = "Sat Feb 17 22:51:35 EST 2024"
$ws.Range("B44").Value = "Sat Feb 17 22:51:47 EST 2024"
$ws.Range("B45").Value = "Sat Feb 17 22:51:59 EST 2024"
$ws.Range("B46").Value = "Sat Feb 17 22:52:12 EST 2024"
$ws.Range("B47").Value = "Sat Feb 17 22:52:24 EST 2024"
$ws.Range("B48").Value = "Sat Feb 17 22:52:36 EST 2024"
$ws.Range("B49").Value = "Sat Feb 17 22:52:48 EST 2024"
$ws.Range("B50").Value = "Sat Feb 17 22:53:00 EST 2024"
$ws.Range("B51").Value = "Sat Feb 17 22:53:12 EST 2024"
$ws.Range("B52").Value = "Sat Feb 17 22:53:25 EST 2024"
$ws.Range("B53").Value = "Sat Feb 17 22:53:38 EST 2024"
$ws.Range("B54").Value = "Sat Feb 17 22:53:50 EST 2024"
